# Scheduled-runner update: refresh market-price-derived profit columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (id 5487)
$ws.Range("H9").Value = 7167.1665
$ws.Range("I9").Value = 8407.267
$ws.Range("J9").Value = 966.6667
$ws.Range("K9").Value = 8407.267
$ws.Range("L9").Value = 966.6667
$ws.Range("M9").Value = -8238.267
$ws.Range("N9").Value = -1304.6667

# Row 31 (id 4576)
$ws.Range("H31").Value = 100329.2
$ws.Range("I31").Value = 125161.5
$ws.Range("K31").Value = 375484.5
$ws.Range("M31").Value = -375254.5

# Row 58 (id 4606)
$ws.Range("H58").Value = 4358
$ws.Range("I58").Value = 428.4
$ws.Range("K58").Value = 1285.2
$ws.Range("M58").Value = -1135.2

# Row 132 (id 44049)
$ws.Range("H132").Value = 2224543.5
$ws.Range("I132").Value = 2858550.2
$ws.Range("K132").Value = 8575650.600000001
$ws.Range("M132").Value = -8573120.600000001

# Row 137 (id 44013)
$ws.Range("H137").Value = 3153.4583
$ws.Range("J137").Value = 4785
$ws.Range("L137").Value = 14355
$ws.Range("N137").Value = -19455

# Row 138 (id 44169)
$ws.Range("H138").Value = 2926.77
$ws.Range("I138").Value = 1424.2115
$ws.Range("J138").Value = 4554.5415
$ws.Range("K138").Value = 4272.6345
$ws.Range("L138").Value = 13663.6245
$ws.Range("M138").Value = 867.3654999999999
$ws.Range("N138").Value = -23943.6245

# Row 141 (id 44161)
$ws.Range("H141").Value = 1172.8085
$ws.Range("I141").Value = 1167.8372
$ws.Range("K141").Value = 3503.5116
$ws.Range("M141").Value = 1676.4884

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (id 27713)
$ws.Range("H2").Value = 1660.4166
$ws.Range("I2").Value = 1592.6
$ws.Range("K2").Value = 1592.6
$ws.Range("M2").Value = -1479.6

# Row 3 (id 2494)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Row 22 (id 2497)
$ws.Range("H22").Value = 10000
$ws.Range("I22").Value = 10000
$ws.Range("K22").Value = 10000
$ws.Range("M22").Value = -9701

# Row 61 (id 43999)
$ws.Range("H61").Value = 2882.762
$ws.Range("I61").Value = 1972.8529
$ws.Range("K61").Value = 1972.8529
$ws.Range("M61").Value = -1760.8529

# Row 116 (id 27713)
$ws.Range("H116").Value = 1660.4166
$ws.Range("I116").Value = 1592.6
$ws.Range("K116").Value = 1592.6
$ws.Range("M116").Value = 701.4000000000001

# Row 122 (id 36168)
$ws.Range("H122").Value = 4591.8125
$ws.Range("I122").Value = 3062.75
$ws.Range("K122").Value = 9188.25
$ws.Range("M122").Value = -6738.25

# Row 132 (id 43997)
$ws.Range("H132").Value = 3733.92
$ws.Range("I132").Value = 1317.4
$ws.Range("K132").Value = 3952.2
$ws.Range("M132").Value = -1422.2

# Row 136 (id 43999)
$ws.Range("H136").Value = 2882.762
$ws.Range("I136").Value = 1972.8529
$ws.Range("K136").Value = 5918.5587
$ws.Range("M136").Value = -3368.5587

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (id 27713)
$ws.Range("H3").Value = 1660.4166
$ws.Range("I3").Value = 1592.6
$ws.Range("K3").Value = 1592.6
$ws.Range("M3").Value = -1478.6

# Row 134 (id 43998)
$ws.Range("H134").Value = 19942.434
$ws.Range("I134").Value = 2150.1064
$ws.Range("K134").Value = 6450.3192
$ws.Range("M134").Value = -3915.3192

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (id 3742)
$ws.Range("H4").Value = 5716500
$ws.Range("J4").Value = 16003200
$ws.Range("L4").Value = 16003200
$ws.Range("N4").Value = -16003424

# Row 23 (id 2703)
$ws.Range("H23").Value = 7375
$ws.Range("I23").Value = 750
$ws.Range("J23").Value = 14000
$ws.Range("K23").Value = 750
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = -510
$ws.Range("N23").Value = -14480

# Row 25 (id 1895)
$ws.Range("H25").Value = 4965.95
$ws.Range("I25").Value = 4962.278
$ws.Range("K25").Value = 4962.278
$ws.Range("M25").Value = -4788.278

# Row 27 (id 2703)
$ws.Range("H27").Value = 7375
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 14000
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = -558
$ws.Range("N27").Value = -14384

# Row 31 (id 44023)
$ws.Range("H31").Value = 177946.17
$ws.Range("J31").Value = 3384.0625
$ws.Range("L31").Value = 3384.0625
$ws.Range("N31").Value = -3974.0625

# Row 34 (id 44023)
$ws.Range("H34").Value = 177946.17
$ws.Range("J34").Value = 3384.0625
$ws.Range("L34").Value = 3384.0625
$ws.Range("N34").Value = -3788.0625

# Row 58 (id 44021)
$ws.Range("H58").Value = 260367.84
$ws.Range("I58").Value = 590449.3
$ws.Range("J58").Value = 5304.909
$ws.Range("K58").Value = 590449.3
$ws.Range("L58").Value = 5304.909
$ws.Range("M58").Value = -590246.3
$ws.Range("N58").Value = -5710.909

# Row 132 (id 44019)
$ws.Range("H132").Value = 4238.394
$ws.Range("I132").Value = 3270.4348
$ws.Range("K132").Value = 9811.304400000001
$ws.Range("M132").Value = -7281.304400000001

# Row 136 (id 44021)
$ws.Range("H136").Value = 260367.84
$ws.Range("I136").Value = 590449.3
$ws.Range("J136").Value = 5304.909
$ws.Range("K136").Value = 1771347.9
$ws.Range("L136").Value = 15914.727
$ws.Range("M136").Value = -1768797.9
$ws.Range("N136").Value = -21014.727

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (id 4847)
$ws.Range("H2").Value = 371.66666
$ws.Range("I2").Value = 260
$ws.Range("K2").Value = 1560
$ws.Range("M2").Value = -1447

# Row 5 (id 43974)
$ws.Range("H5").Value = 550.7222
$ws.Range("I5").Value = 350.875
$ws.Range("J5").Value = 2149.5
$ws.Range("K5").Value = 1052.625
$ws.Range("L5").Value = 6448.5
$ws.Range("M5").Value = -940.625
$ws.Range("N5").Value = -6672.5

# Row 17 (id 4640)
$ws.Range("H17").Value = 1621.3
$ws.Range("I17").Value = 1589.125
$ws.Range("K17").Value = 4767.375
$ws.Range("M17").Value = -4598.375

# Row 75 (id 12863)
$ws.Range("H75").Value = 1407.6666
$ws.Range("I75").Value = 1048.3334
$ws.Range("J75").Value = 1767
$ws.Range("K75").Value = 3145.0002
$ws.Range("L75").Value = 5301
$ws.Range("M75").Value = -2147.0002
$ws.Range("N75").Value = -7297

# Row 78 (id 12863)
$ws.Range("H78").Value = 1407.6666
$ws.Range("I78").Value = 1048.3334
$ws.Range("J78").Value = 1767
$ws.Range("K78").Value = 9435.000599999999
$ws.Range("L78").Value = 15903
$ws.Range("M78").Value = -4443.000599999999
$ws.Range("N78").Value = -25887

# Row 135 (id 43974)
$ws.Range("H135").Value = 550.7222
$ws.Range("I135").Value = 350.875
$ws.Range("J135").Value = 2149.5
$ws.Range("K135").Value = 3157.875
$ws.Range("L135").Value = 19345.5
$ws.Range("M135").Value = -622.875
$ws.Range("N135").Value = -24415.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (id 5062)
$ws.Range("H2").Value = 157.55556
$ws.Range("I2").Value = 224.75
$ws.Range("K2").Value = 224.75
$ws.Range("M2").Value = -111.75

# Row 102 (id 36169)
$ws.Range("H102").Value = 2236.0454
$ws.Range("I102").Value = 1513.7273
$ws.Range("J102").Value = 2958.3635
$ws.Range("K102").Value = 1513.7273
$ws.Range("L102").Value = 2958.3635
$ws.Range("M102").Value = 108.2727
$ws.Range("N102").Value = -6202.363499999999

# Row 132 (id 44008)
$ws.Range("H132").Value = 528713.3
$ws.Range("I132").Value = 774121.4
$ws.Range("K132").Value = 2322364.2
$ws.Range("M132").Value = -2319834.2

$ws = $wb.Worksheets.Item("LTW")
# Row 58 (id 1728)
$ws.Range("H58").Value = 20066.334
$ws.Range("I58").Value = 19100
$ws.Range("K58").Value = 19100
$ws.Range("M58").Value = -18840

# Row 132 (id 44058)
$ws.Range("H132").Value = 3226.2942
$ws.Range("I132").Value = 2811.5715
$ws.Range("K132").Value = 8434.7145
$ws.Range("M132").Value = -5904.7145

# Row 135 (id 42036)
$ws.Range("H135").Value = 57149.4
$ws.Range("J135").Value = 57149.4
$ws.Range("L135").Value = 57149.4
$ws.Range("N135").Value = -67289.39999999999

# Row 136 (id 44060)
$ws.Range("H136").Value = 6248.9165
$ws.Range("I136").Value = 5754.625
$ws.Range("J136").Value = 7237.5
$ws.Range("K136").Value = 17263.875
$ws.Range("L136").Value = 21712.5
$ws.Range("M136").Value = -14713.875
$ws.Range("N136").Value = -26812.5

$ws = $wb.Worksheets.Item("WVR")
# Row 93 (id 19613)
$ws.Range("H93").Value = 43750
$ws.Range("J93").Value = 43750
$ws.Range("L93").Value = 43750
$ws.Range("N93").Value = -48742

# Row 122 (id 36208)
$ws.Range("H122").Value = 38465880
$ws.Range("I122").Value = 55558716
$ws.Range("K122").Value = 166676148
$ws.Range("M122").Value = -166673698

# Row 132 (id 44029)
$ws.Range("H132").Value = 26401.414
$ws.Range("J132").Value = 104879.9
$ws.Range("L132").Value = 314639.7
$ws.Range("N132").Value = -319699.7

# Row 136 (id 44031)
$ws.Range("H136").Value = 93292.82000000001
$ws.Range("I136").Value = 1956.9231
$ws.Range("K136").Value = 5870.7693
$ws.Range("M136").Value = -3320.7693
